$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 211, shifting existing rows 211:330 down to 212:331
$ws.Rows("211:211").Insert()

# Copy the date cell style (style index 2 / numFmt 165) from the row below (now row 212)
$ws.Range("D212").Copy()
$ws.Range("D211").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row 211 with data
$ws.Range("A211").Value = 6
$ws.Range("B211").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C211").Value = "Metropolitana"
$ws.Range("D211").Value = 44438
$ws.Range("E211").Value = 13
$ws.Range("F211").Value = 100112003
$ws.Range("G211").Value = "Ajo"
$ws.Range("H211").Value = "Chino"
$ws.Range("I211").Value = "Primera"
$ws.Range("J211").Value = 3700
$ws.Range("K211").Value = 14500
$ws.Range("L211").Value = 15000
$ws.Range("M211").Value = 14703
$ws.Range("N211").Value = "`$/malla 10 kilos"
$ws.Range("O211").Value = "China"
$ws.Range("P211").Value = 1470
$ws.Range("Q211").Value = 10
$ws.Range("R211").Value = "Hortaliza"
